# Update a few numeric "views/popularity" counts (column F) on two sheets:
# "展览" (Worksheets index 1) and "全部类型" (Worksheets index 4).
# The same events appear on both sheets (the latter is an "all types" roll-up),
# so both copies of each value need to be bumped by the same amount.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6605
$wsExhibit.Range("F13").Value = 1275
$wsExhibit.Range("F14").Value = 6
$wsExhibit.Range("F15").Value = 3284
$wsExhibit.Range("F18").Value = 1922
$wsExhibit.Range("F19").Value = 36
$wsExhibit.Range("F20").Value = 26

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6605
$wsAll.Range("F14").Value = 1275
$wsAll.Range("F15").Value = 6
$wsAll.Range("F16").Value = 3284
$wsAll.Range("F19").Value = 1922
$wsAll.Range("F20").Value = 36
$wsAll.Range("F21").Value = 26
